# From v1.2.3 to v1.2.4
# The TC3 and TC4 test-case scenarios are swapped:
#  - TC3 ("cancelar diária") becomes the "detalhar diária" scenario
#  - TC4 ("detalhar diária") becomes the "cancelar diária" scenario
# Test-case ID labels (TC3 in row 21 block, TC4 in row 28 block) stay put;
# only the Steps / Expected Results text is exchanged between the two blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc3Steps    = $ws.Range("B25").Value2
$tc3Expected = $ws.Range("D25").Value2
$tc4Steps    = $ws.Range("B32").Value2
$tc4Expected = $ws.Range("D32").Value2

$ws.Range("B25").Value2 = $tc4Steps
$ws.Range("D25").Value2 = $tc4Expected
$ws.Range("B32").Value2 = $tc3Steps
$ws.Range("D32").Value2 = $tc3Expected
